$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2.05
$ws.Range("S2").Value = 2.35
$ws.Range("T2").Value = 1.57

# Row 4
$ws.Range("G4").Value = 1.45
$ws.Range("J4").Value = 2.05
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 3
$ws.Range("S4").Value = 2.2
$ws.Range("T4").Value = 1.67
$ws.Range("AN4").Value = 15
$ws.Range("AP4").Value = 26

# Row 5
$ws.Range("G5").Value = 2.55
$ws.Range("I5").Value = 2.9
$ws.Range("J5").Value = 3.4
$ws.Range("AF5").Value = 26
$ws.Range("AG5").Value = 23
$ws.Range("AH5").Value = 41
$ws.Range("AI5").Value = 7
$ws.Range("AM5").Value = 451
$ws.Range("AN5").Value = 7.5
$ws.Range("AP5").Value = 11
$ws.Range("AQ5").Value = 29
$ws.Range("AR5").Value = 26

# Row 6
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("S6").Value = 1.85
$ws.Range("T6").Value = 2

# Row 7
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
$ws.Range("S7").Value = 1.85
$ws.Range("T7").Value = 2
